# Append 9 new vocabulary entries to the "words" sheet (rows 63-71),
# all logged on 2020-12-03 - the next day's batch following the existing
# daily study-log pattern already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "관습",
    "긴밀하다",
    "지도자",
    "거듭",
    "강조하다",
    "설립하다",
    "몰수되다",
    "영지",
    "봉기"
)

$meanings = @(
    "1.) n. custom, convention",
    "1.) adj. close, intimate; 2.) adj. tight, close",
    "1.) n. leader, guide",
    "1.) adj. again, once more, once again, repeatedly",
    "1.) v. emphasize, stress",
    "1.) v. establish, found",
    "1.) v. be confiscated, be forfeited, be sequestered",
    "1.) n. territory, possession, dominion; 2.) fief, feud, vassalage, estate",
    "1.) n. uprising, revolt, rebellion"
)

$newDate = "2020-12-03"
$startRow = 63
$count = $words.Length

# Write column A (FOREIGN / Korean word) for every new row first, then
# column B (ENGLISH meaning), then column C (DATE) - matches the grouping
# order the source data was appended in.
for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $words[$i]
}

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $ws.Range("B$row").Value = $meanings[$i]
}

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    # Force text format before assigning so the engine doesn't coerce the
    # "yyyy-mm-dd"-shaped string into a date serial number (column C is
    # plain text everywhere else in the sheet).
    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = $newDate
    $ws.Range("C$row").Style = "Normal"
}
